$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Job")

# Rename the "是否常规职业" column header to "是否特殊职业".
# (Built from code points so the literal Chinese text survives the PS host's
# string handling; this also renames the backing Excel Table column since the
# header cell and the ListObject column name stay in sync.)
$newHeader = [string]([char]0x662F) + [string]([char]0x5426) + [string]([char]0x7279) + [string]([char]0x6B8A) + [string]([char]0x804C) + [string]([char]0x4E1A)
$ws.Range("L1").Value = $newHeader

# Mark rows 9, 13 and 14 (job ids 11000005 / 11000009 / 11000010) as special
# jobs too, matching row 4 / 15 / 16 which already say "true". Copy/paste
# special (values only) from L4 so the cell keeps its literal text "true"
# instead of Excel auto-converting the typed word into a Boolean.
$ws.Range("L4").Copy()
$ws.Range("L9").PasteSpecial(-4163)
$ws.Range("L13").PasteSpecial(-4163)
$ws.Range("L14").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Move the active selection on the Job sheet to L14.
$ws.Activate()
$ws.Range("L14").Select()
